# Rename quizzes to meaningful names; update metadata sheet and fix a
# mis-tagged subject; move the "active tab" selection from quiz6 (now
# quiz_politics2) to metadata_quiz.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the quiz worksheets (tab names) ---------------------------
$wb.Worksheets.Item("quiz1").Name = "quiz_geography1"
$wb.Worksheets.Item("quiz2").Name = "quiz_politics1"
$wb.Worksheets.Item("quiz3").Name = "quiz_environment1"
$wb.Worksheets.Item("quiz4").Name = "quiz_independence1"
$wb.Worksheets.Item("quiz5").Name = "quiz_environment2"
$wb.Worksheets.Item("quiz6").Name = "quiz_politics2"

# --- 2. Update the metadata_quiz sheet so its quiz-name column matches ---
$ws = $wb.Worksheets.Item("metadata_quiz")

$ws.Range("A2").Value = "quiz_geography1"
$ws.Range("A3").Value = "quiz_politics1"
$ws.Range("A4").Value = "quiz_environment1"
$ws.Range("A5").Value = "quiz_independence1"
$ws.Range("A6").Value = "quiz_environment2"
$ws.Range("A7").Value = "quiz_politics2"

# Fix subject mismatch: quiz_environment1 was tagged INTA, should be ENVI
$ws.Range("C4").Value = "ENVI"

# --- 3. Re-point the active/selected tab + column width on metadata_quiz -
$ws.Columns.Item(1).AutoFit()
$ws.Range("A2").Select()
$ws.Activate()

# --- 4. quiz_politics2 (former quiz6) is no longer the selected tab ------
$wsPol2 = $wb.Worksheets.Item("quiz_politics2")
$wsPol2.Columns.Item(1).AutoFit()
$wsPol2.Columns.Item(2).AutoFit()
$wsPol2.Columns.Item(3).AutoFit()
$wsPol2.Columns.Item(12).AutoFit()
$wsPol2.Columns.Item(13).AutoFit()

# Re-activate metadata_quiz last so it ends up as the workbook's active tab
$ws.Activate()
$ws.Range("A2").Select()
